# Update the embedded build timestamp throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$oldVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $oldTimestamp)"
$newVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersionString"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Shoal Creek Coal Mine, United States, M1068, version '$newVersionString'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds the build_version value per-row (header in S1, data rows S2:S26).
$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value() -eq $oldVersionString) {
        $cell.Value = $newVersionString
    }
}
